$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-run TPM pipeline output for S100a9-Tlr4 ---
# Ligand-side measurements (E:J) depend only on the "Sending cluster" (col A),
# so they are identical across every block of rows that share a sending cluster.

# Sending cluster = ECs (rows 2-6)
$ws.Range("E2:E6").Value = 1
$ws.Range("F2:F6").Value = 0.3333333333333333
$ws.Range("G2:G6").Value = 0.7540426666666667
$ws.Range("H2:H6").Value = 2.262128
$ws.Range("I2:I6").Value = 0.06765183538434597
$ws.Range("J2:J6").Value = 0.06765183538434595

# Sending cluster = Inflammatory-Mac (rows 7-11)
$ws.Range("G7:G11").Value = 8.291411999999999
$ws.Range("H7:H11").Value = 24.874236
$ws.Range("I7:I11").Value = 0.7438958888194531
$ws.Range("J7:J11").Value = 0.7438958888194531

# Sending cluster = Resolving-Mac (rows 12-16)
$ws.Range("G12:G16").Value = 2.100476
$ws.Range("H12:H16").Value = 6.301428
$ws.Range("I12:I16").Value = 0.1884522757962009
$ws.Range("J12:J16").Value = 0.1884522757962009

# Receptor-side measurements (M:P) depend only on the "Target cluster" (col D),
# and the edge weights/specificities (Q:T) are unique per row.

# Row 2
$ws.Range("M2").Value = 7.658574666666667
$ws.Range("N2").Value = 22.975724
$ws.Range("O2").Value = 0.1056247585863608
$ws.Range("P2").Value = 0.1080674328374073
$ws.Range("Q2").Value = 5.774892064519111
$ws.Range("R2").Value = 51.97402858067201
$ws.Range("S2").Value = 0.007145708780395766
$ws.Range("T2").Value = 0.007310960176725141

# Row 3
$ws.Range("O3").Value = 0.1485332542774742
$ws.Range("P3").Value = 0.1519682288090475
$ws.Range("Q3").Value = 8.120856538979556
$ws.Range("R3").Value = 73.087708850816
$ws.Range("S3").Value = 0.01004854726748089
$ws.Range("T3").Value = 0.0102809295990403

# Row 4
$ws.Range("M4").Value = 21.58649266666667
$ws.Range("N4").Value = 64.759478
$ws.Range("O4").Value = 0.2977144150029286
$ws.Range("P4").Value = 0.3045993475265701
$ws.Range("Q4").Value = 16.27713649435378
$ws.Range("R4").Value = 146.494228449184
$ws.Range("S4").Value = 0.02014092659532498
$ws.Range("T4").Value = 0.0206067049170467

# Row 5
$ws.Range("M5").Value = 4.9167055
$ws.Range("N5").Value = 9.833411
$ws.Range("O5").Value = 0.06780972362103574
$ws.Range("P5").Value = 0.04625192585030714
$ws.Range("Q5").Value = 3.707405726434667
$ws.Range("R5").Value = 22.244434358608
$ws.Range("S5").Value = 0.004587452259868306
$ws.Range("T5").Value = 0.003129027673833954

# Row 6
$ws.Range("M6").Value = 27.57585133333333
$ws.Range("N6").Value = 82.727554
$ws.Range("O6").Value = 0.3803178485122005
$ws.Range("P6").Value = 0.3891130649766679
$ws.Range("Q6").Value = 20.79336847499022
$ws.Range("R6").Value = 187.140316274912
$ws.Range("S6").Value = 0.02572920048127602
$ws.Range("T6").Value = 0.02632421301769984

# Row 7
$ws.Range("M7").Value = 7.658574666666667
$ws.Range("N7").Value = 22.975724
$ws.Range("O7").Value = 0.1056247585863608
$ws.Range("P7").Value = 0.1080674328374073
$ws.Range("Q7").Value = 63.500397894096
$ws.Range("R7").Value = 571.5035810468639
$ws.Range("S7").Value = 0.07857382366994106
$ws.Range("T7").Value = 0.08039091900301967

# Row 8
$ws.Range("O8").Value = 0.1485332542774742
$ws.Range("P8").Value = 0.1519682288090475
$ws.Range("Q8").Value = 89.29649519068799
$ws.Range("R8").Value = 803.6684567161919
$ws.Range("S8").Value = 0.1104932772099875
$ws.Range("T8").Value = 0.1130485406422244

# Row 9
$ws.Range("M9").Value = 21.58649266666667
$ws.Range("N9").Value = 64.759478
$ws.Range("O9").Value = 0.2977144150029286
$ws.Range("P9").Value = 0.3045993475265701
$ws.Range("Q9").Value = 178.982504334312
$ws.Range("R9").Value = 1610.842539008808
$ws.Range("S9").Value = 0.2214685293629671
$ws.Range("T9").Value = 0.2265902023621033

# Row 10
$ws.Range("M10").Value = 4.9167055
$ws.Range("N10").Value = 9.833411
$ws.Range("O10").Value = 0.06780972362103574
$ws.Range("P10").Value = 0.04625192585030714
$ws.Range("Q10").Value = 40.766430983166
$ws.Range("R10").Value = 244.598585898996
$ws.Range("S10").Value = 0.05044337462367185
$ws.Range("T10").Value = 0.03440661749002567

# Row 11
$ws.Range("M11").Value = 27.57585133333333
$ws.Range("N11").Value = 82.727554
$ws.Range("O11").Value = 0.3803178485122005
$ws.Range("P11").Value = 0.3891130649766679
$ws.Range("Q11").Value = 228.642744655416
$ws.Range("R11").Value = 2057.784701898744
$ws.Range("S11").Value = 0.2829168839528856
$ws.Range("T11").Value = 0.28945960932208

# Row 12
$ws.Range("M12").Value = 7.658574666666667
$ws.Range("N12").Value = 22.975724
$ws.Range("O12").Value = 0.1056247585863608
$ws.Range("P12").Value = 0.1080674328374073
$ws.Range("Q12").Value = 16.08665228154133
$ws.Range("R12").Value = 144.779870533872
$ws.Range("S12").Value = 0.01990522613602401
$ws.Range("T12").Value = 0.0203655536576625

# Row 13
$ws.Range("O13").Value = 0.1485332542774742
$ws.Range("P13").Value = 0.1519682288090475
$ws.Range("Q13").Value = 22.62161680449066
$ws.Range("R13").Value = 203.594551240416
$ws.Range("S13").Value = 0.02799142980000581
$ws.Range("T13").Value = 0.02863875856778279

# Row 14
$ws.Range("M14").Value = 21.58649266666667
$ws.Range("N14").Value = 64.759478
$ws.Range("O14").Value = 0.2977144150029286
$ws.Range("P14").Value = 0.3045993475265701
$ws.Range("Q14").Value = 45.34190977050934
$ws.Range("R14").Value = 408.077187934584
$ws.Range("S14").Value = 0.0561049590446365
$ws.Range("T14").Value = 0.05740244024742003

# Row 15
$ws.Range("M15").Value = 4.9167055
$ws.Range("N15").Value = 9.833411
$ws.Range("O15").Value = 0.06780972362103574
$ws.Range("P15").Value = 0.04625192585030714
$ws.Range("Q15").Value = 10.327421901818
$ws.Range("R15").Value = 61.96453141090799
$ws.Range("S15").Value = 0.01277889673749559
$ws.Range("T15").Value = 0.008716280686447515

# Row 16
$ws.Range("M16").Value = 27.57585133333333
$ws.Range("N16").Value = 82.727554
$ws.Range("O16").Value = 0.3803178485122005
$ws.Range("P16").Value = 0.3891130649766679
$ws.Range("Q16").Value = 57.92241390523466
$ws.Range("R16").Value = 521.3017251471119
$ws.Range("S16").Value = 0.07167176407803896
$ws.Range("T16").Value = 0.07332924263688805
